$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2 through 481
# from serial date 45202 (2023-10-03) to 45203 (2023-10-04)
$ws.Range("C2:C481").Value = 45203
